$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(65, 45511, 6.1179699999999997, 5.6033099999999996, "NA", "NA", 6.12, 5.61),
    @(66, 45511, 6.1429999999999998, 5.6236600000000001, "NA", "NA", 6.15, 5.63),
    @(67, 45511, 6.1602499999999996, 5.6373899999999999, "NA", "NA", 6.17, 5.64),
    @(68, 45512, 6.1177400000000004, 5.617, "NA", "NA", 6.12, 5.62),
    @(69, 45512, 6.0851300000000004, 5.5742599999999998, "NA", "NA", 6.09, 5.58),
    @(70, 45512, 6.0570500000000003, 5.5475099999999999, "NA", "NA", 6.06, 5.55),
    @(71, 45513, 6.0171000000000001, 5.5172400000000001, "NA", "NA", 6.02, 5.52),
    @(72, 45513, 6.0211600000000001, 5.5146600000000001, "NA", "NA", 6.03, 5.52),
    @(73, 45513, 6.0135500000000004, 5.5087000000000002, 6.14628, 5.6179699999999997, 6.03, 5.51),
    @(74, 45514, 6.0135500000000004, 5.5087000000000002, 6.14628, 5.6179699999999997, 6.03, 5.51),
    @(75, 45514, 6.0135500000000004, 5.5087000000000002, 6.14628, 5.6179699999999997, 6.02, 5.51),
    @(76, 45514, 6.0135500000000004, 5.5087000000000002, 6.14628, 5.6179699999999997, 6.02, 5.51),
    @(77, 45515, 6.0078100000000001, 5.4993999999999996, 6.1425000000000001, 5.62113, 6.02, 5.5),
    @(78, 45515, 6.0093300000000003, 5.4977400000000003, 6.1349600000000004, 5.6053800000000003, 6.01, 5.49),
    @(79, 45515, 6.0075200000000004, 5.4938399999999996, 6.1652199999999997, 5.6053800000000003, 6.02, 5.49)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}

# Match the author's final view/selection state after appending the new rows.
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 55 } catch {}
$ws.Range("D84").Select() | Out-Null

